$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Meetings » NeoGrowth Credit Pvt.Ltd."
$ws.Range("A4").Value = "Meeting Page "

$ws.Range("A7").Select()
